$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change the "coordonne y" value used for the xytheta computation from 90 to 0.
# This makes L17 (=K17/H17) and M17 (=(1-L17)*100) evaluate to #DIV/0!
# and makes H23 (=H17-H21) evaluate to -45 instead of 45.
$ws.Range("H17").Value = 0

# Force Excel to recalculate all formulas so the dependent cells (L17, M17, H23)
# pick up the new #DIV/0! / -45 results.
$excel.CalculateFullRebuild()

# Update the view: scroll so row 3 is the top-left visible row and move the
# active selection to H19.
$window = $ws.Application.ActiveWindow
$window.ScrollRow = 3
$ws.Range("H19").Select()
